$wb = $excel.ActiveWorkbook

# --- Sheet "Tests" (sheet1) ---
$ws1 = $wb.Worksheets.Item(1)

# Swap the "execute" flag between row 2 and row 6
$ws1.Range("C2").Value = "no"
$ws1.Range("C6").Value = "yes"

# Move the active selection to C6
$ws1.Range("C6").Select()

# --- Sheet "DataProviderTests" (sheet2) ---
$ws2 = $wb.Worksheets.Item(2)

# Flip a few existing "execute" flags
$ws2.Range("B2").Value = "no"
$ws2.Range("B3").Value = "no"
$ws2.Range("B10").Value = "yes"

# Add two new rows (11 and 12) mirroring the style of row 10
$ws2.Range("A10:E10").Copy($ws2.Range("A11:E11"))
$ws2.Range("A10:E10").Copy($ws2.Range("A12:E12"))

$ws2.Range("A11").Value = "navigateToTabletsPageViaHamburgerMenuTest"
$ws2.Range("B11").Value2 = "yes"
$ws2.Range("C11").Value = "firefox"
$ws2.Range("D11").Value2 = "'"
$ws2.Range("E11").Value2 = "'"

$ws2.Range("A12").Value = "navigateToTabletsPageViaHamburgerMenuTest"
$ws2.Range("B12").Value2 = "yes"
$ws2.Range("C12").Value = "edge"
$ws2.Range("D12").Value2 = "'"
$ws2.Range("E12").Value2 = "'"

# Update the used-range dimension and active selection
$ws2.Range("C12").Select()
